$wb = $excel.ActiveWorkbook

# --- Page setup (paper size / orientation) for the first three sheets ---
# "Variation eta", "Evolution sigma", "Evolution exponentielle"
for ($i = 1; $i -le 3; $i++) {
    $s = $wb.Worksheets.Item($i)
    $s.PageSetup.PaperSize = 9
    $s.PageSetup.Orientation = 1
}

# --- Sheet "Evolution distribution" (5th sheet): add the 3 new data rows ---
$ws5 = $wb.Worksheets.Item(5)

# Write B8 first, then B6, then B7 so that new shared-string entries are
# created in the same order as in the target workbook (indices 17,18,19).
$ws5.Range("B8").Value = "[−1, 1] × [−1, 1] et 3* [0, 1] x [0, 1]"
$ws5.Range("B6").Value = "[−1, 0] × [-1, 1] et 3*[0, 1] × [-1, 1] "
$ws5.Range("B7").Value = "[−1, 0] × [-1, 1] et 7*[0, 1] × [-1, 1] "

$ws5.Range("C6").Value = 0.013693043211208801
$ws5.Range("D6").Value = 5.5910500134726098

$ws5.Range("C7").Value = 0.0161745947705098
$ws5.Range("D7").Value = 5.1341867715500697

$ws5.Range("C8").Value = 0.015226258240725101
$ws5.Range("D8").Value = 4.6756390201862503

# Copy the formatting of the row above onto the new rows so the new cells
# keep the same (centered) cell style as the rest of the table.
$ws5.Range("B5:D5").Copy()
$ws5.Range("B6:D8").PasteSpecial(-4122)

# Update the selection on this sheet and make it the active one, mirroring
# the author's final view/state.
$ws5.Range("D10").Select()
$ws5.Activate()
